# Generate Report for Handback
# Updates handback status timestamps / priority values produced by a new
# report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-28 06:19:35"
$overview.Range("G4").Value = "2016-08-28 06:19:35"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-28 06:19:31"
$zhcn.Range("H4").Value = "2016-08-28 06:19:31"
$zhcn.Range("K3").Value = "2016-08-28 06:19:48"
$zhcn.Range("K4").Value = "2016-08-28 06:19:48"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "mt"
$dede.Range("E4").Value = "mt"
$dede.Range("H3").Value = "2016-08-28 06:19:35"
$dede.Range("H4").Value = "2016-08-28 06:19:35"
$dede.Range("K3").Value = "2016-08-28 06:19:54"
$dede.Range("K4").Value = "2016-08-28 06:19:54"
